$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''76.350.33'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.37%  '

# Row 3
$ws.Range('D3').Value = '''3.038.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.72%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = '''200.53'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.30%  '

# Row 6
$ws.Range('D6').Value = '''623.16'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.29%  '

# Row 7
$ws.Range('D7').Value = '''0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').Value = '''0.550'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.07%  '

# Row 9
$ws.Range('D9').Value = '''0.205'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.86%  '

# Row 10
$ws.Range('D10').Value = '''3.036.86'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.69%  '

# Row 11
$ws.Range('D11').Value = '''0.440'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.72%  '

# Row 12
$ws.Range('D12').Value = '''0.160'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.70%  '

# Row 13
$ws.Range('D13').Value = '''5.24'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.97%  '

# Row 14
$ws.Range('D14').Value = '''3.599.32'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.75%  '

# Row 15
$ws.Range('D15').Value = '''29.19'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.92%  '

# Row 16
$ws.Range('D16').Value = '''76.311.69'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.45%  '

# Row 17
$ws.Range('D17').Value = '''0.0000192'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.59%  '

# Row 18
$ws.Range('D18').Value = '''3.050.67'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.46%  '

# Row 19
$ws.Range('D19').Value = '''13.54'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.89%  '

# Row 20
$ws.Range('D20').Value = '''8.99'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.66%  '

# Row 21
$ws.Range('D21').Value = '''375.01'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.41%  '

# Row 22
$ws.Range('D22').Value = '''2.30'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.42%  '

# Row 23
$ws.Range('D23').Value = '''4.35'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.92%  '

# Row 24
$ws.Range('D24').Value = '''73.49'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.69%  '

# Row 25
$ws.Range('D25').Value = '''3.200.07'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.82%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''0.997'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.24%  '

# Row 27
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = '''4.37'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.02%  '

# Row 28
$ws.Range('D28').Value = '''9.81'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.48%  '

# Row 29
$ws.Range('D29').Value = '''0.0000110'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.27%  '

# Row 30
$ws.Range('D30').Value = '''1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.01%  '

# Row 31
$ws.Range('D31').Value = '''8.24'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.29%  '

# Row 32
$ws.Range('D32').Value = '''1.41'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.46%  '

# Row 33
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''1.95'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +6.01%  '

# Row 34
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '''496.42'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.26%  '

# Row 35
$ws.Range('D35').Value = '''0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.03%  '

# Row 36
$ws.Range('D36').Value = '''20.63'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.73%  '

# Row 37
$ws.Range('D37').Value = '''162.84'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.60%  '

# Row 38
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').Value = '''20.04'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.17%  '

# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '''0.384'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.69%  '

# Row 40
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.116'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.65%  '

# Row 41
$ws.Range('D41').Value = '''189.60'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.08%  '

# Row 42
$ws.Range('D42').Value = '''0.105'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.53%  '

# Row 44
$ws.Range('D44').Value = '''0.797'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +20.64%  '

# Row 45
$ws.Range('D45').Value = '''5.11'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.88%  '

# Row 46
$ws.Range('D46').Value = '''1.27'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.64%  '

# Row 47
$ws.Range('D47').Value = '''42.01'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.38%  '

# Row 48
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '''1.65'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.83%  '

# Row 49
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '''2.47'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.38%  '

# Row 50
$ws.Range('D50').Value = '''0.608'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.94%  '

# Row 51
$ws.Range('D51').Value = '''3.90'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.06%  '
